$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 1040.682692307692
$ws.Cells.Item(2, 2).Value = 1054.632075471698
$ws.Cells.Item(2, 3).Value = 1054.461538461539
$ws.Cells.Item(2, 4).Value = 1052.326923076923
$ws.Cells.Item(2, 5).Value = 1058.807692307692
$ws.Cells.Item(2, 6).Value = 1049.951923076923
$ws.Cells.Item(2, 7).Value = 1043.423076923077
$ws.Cells.Item(3, 1).Value = 1037.105769230769
$ws.Cells.Item(3, 2).Value = 1054.179245283019
$ws.Cells.Item(3, 3).Value = 1051.480769230769
$ws.Cells.Item(3, 4).Value = 1051.653846153846
$ws.Cells.Item(3, 5).Value = 1057.048076923077
$ws.Cells.Item(3, 6).Value = 1048.288461538461
$ws.Cells.Item(3, 7).Value = 1043.475961538461
$ws.Cells.Item(4, 1).Value = 1037.009615384615
$ws.Cells.Item(4, 2).Value = 1051.783018867925
$ws.Cells.Item(4, 3).Value = 1049.826923076923
$ws.Cells.Item(4, 4).Value = 1050.432692307692
$ws.Cells.Item(4, 5).Value = 1054.548076923077
$ws.Cells.Item(4, 6).Value = 1046.307692307692
$ws.Cells.Item(4, 7).Value = 1040.813725490196
$ws.Cells.Item(5, 1).Value = 1038.038461538461
$ws.Cells.Item(5, 2).Value = 1054.915094339623
$ws.Cells.Item(5, 3).Value = 1051.201923076923
$ws.Cells.Item(5, 4).Value = 1049.961538461539
$ws.Cells.Item(5, 5).Value = 1055.538461538461
$ws.Cells.Item(5, 6).Value = 1045.134615384615
$ws.Cells.Item(5, 7).Value = 1042.125
$ws.Cells.Item(6, 1).Value = 1037.317307692308
$ws.Cells.Item(6, 2).Value = 1052.207547169811
$ws.Cells.Item(6, 3).Value = 1050.471153846154
$ws.Cells.Item(6, 4).Value = 1047.307692307692
$ws.Cells.Item(6, 5).Value = 1052.509615384615
$ws.Cells.Item(6, 6).Value = 1044
$ws.Cells.Item(6, 7).Value = 1041.894230769231
$ws.Cells.Item(7, 1).Value = 1038.115384615385
$ws.Cells.Item(7, 2).Value = 1055.38679245283
$ws.Cells.Item(7, 3).Value = 1049.971153846154
$ws.Cells.Item(7, 4).Value = 1055.375
$ws.Cells.Item(7, 5).Value = 1049.403846153846
$ws.Cells.Item(7, 6).Value = 1041.807692307692
$ws.Cells.Item(7, 7).Value = 1043.971153846154
$ws.Cells.Item(8, 1).Value = 1051.625
$ws.Cells.Item(8, 2).Value = 1066.698113207547
$ws.Cells.Item(8, 3).Value = 1060.307692307692
$ws.Cells.Item(8, 4).Value = 1063.519230769231
$ws.Cells.Item(8, 5).Value = 1060.394230769231
$ws.Cells.Item(8, 6).Value = 1048.317307692308
$ws.Cells.Item(8, 7).Value = 1051.153846153846
$ws.Cells.Item(9, 1).Value = 1073.644230769231
$ws.Cells.Item(9, 2).Value = 1095.047169811321
$ws.Cells.Item(9, 3).Value = 1079.25
$ws.Cells.Item(9, 4).Value = 1086.971153846154
$ws.Cells.Item(9, 5).Value = 1082.192307692308
$ws.Cells.Item(9, 6).Value = 1055.778846153846
$ws.Cells.Item(9, 7).Value = 1072.913461538461
$ws.Cells.Item(10, 1).Value = 1143.25
$ws.Cells.Item(10, 2).Value = 1147.377358490566
$ws.Cells.Item(10, 3).Value = 1174.980769230769
$ws.Cells.Item(10, 4).Value = 1161.009615384615
$ws.Cells.Item(10, 5).Value = 1138.798076923077
$ws.Cells.Item(10, 6).Value = 1064.730769230769
$ws.Cells.Item(10, 7).Value = 1081.884615384615
$ws.Cells.Item(11, 1).Value = 1213.605769230769
$ws.Cells.Item(11, 2).Value = 1230.924528301887
$ws.Cells.Item(11, 3).Value = 1242.461538461539
$ws.Cells.Item(11, 4).Value = 1224.625
$ws.Cells.Item(11, 5).Value = 1190.067307692308
$ws.Cells.Item(11, 6).Value = 1087.355769230769
$ws.Cells.Item(11, 7).Value = 1103.932692307692
$ws.Cells.Item(12, 1).Value = 1247.769230769231
$ws.Cells.Item(12, 2).Value = 1281.481132075472
$ws.Cells.Item(12, 3).Value = 1282.759615384615
$ws.Cells.Item(12, 4).Value = 1267.423076923077
$ws.Cells.Item(12, 5).Value = 1216.461538461539
$ws.Cells.Item(12, 6).Value = 1104.403846153846
$ws.Cells.Item(12, 7).Value = 1128.442307692308
$ws.Cells.Item(13, 1).Value = 1264.75
$ws.Cells.Item(13, 2).Value = 1295.584905660377
$ws.Cells.Item(13, 3).Value = 1304.596153846154
$ws.Cells.Item(13, 4).Value = 1278.105769230769
$ws.Cells.Item(13, 5).Value = 1242.461538461539
$ws.Cells.Item(13, 6).Value = 1104.288461538461
$ws.Cells.Item(13, 7).Value = 1116.009615384615
$ws.Cells.Item(14, 1).Value = 1278.086538461539
$ws.Cells.Item(14, 2).Value = 1304.349056603774
$ws.Cells.Item(14, 3).Value = 1305.894230769231
$ws.Cells.Item(14, 4).Value = 1279.634615384615
$ws.Cells.Item(14, 5).Value = 1234.778846153846
$ws.Cells.Item(14, 6).Value = 1129.307692307692
$ws.Cells.Item(14, 7).Value = 1123.788461538461
$ws.Cells.Item(15, 1).Value = 1278.153846153846
$ws.Cells.Item(15, 2).Value = 1318.216981132075
$ws.Cells.Item(15, 3).Value = 1319.528846153846
$ws.Cells.Item(15, 4).Value = 1283.961538461539
$ws.Cells.Item(15, 5).Value = 1228.990384615385
$ws.Cells.Item(15, 6).Value = 1116.894230769231
$ws.Cells.Item(15, 7).Value = 1128.682692307692
$ws.Cells.Item(16, 1).Value = 1271.288461538461
$ws.Cells.Item(16, 2).Value = 1314.216981132075
$ws.Cells.Item(16, 3).Value = 1304.721153846154
$ws.Cells.Item(16, 4).Value = 1273.153846153846
$ws.Cells.Item(16, 5).Value = 1194.163461538461
$ws.Cells.Item(16, 6).Value = 1085.769230769231
$ws.Cells.Item(16, 7).Value = 1107.182692307692
$ws.Cells.Item(17, 1).Value = 1268.769230769231
$ws.Cells.Item(17, 2).Value = 1278.641509433962
$ws.Cells.Item(17, 3).Value = 1261.115384615385
$ws.Cells.Item(17, 4).Value = 1238.009615384615
$ws.Cells.Item(17, 5).Value = 1179.875
$ws.Cells.Item(17, 6).Value = 1075.836538461539
$ws.Cells.Item(17, 7).Value = 1086.346153846154
$ws.Cells.Item(18, 1).Value = 1196.769230769231
$ws.Cells.Item(18, 2).Value = 1201.688679245283
$ws.Cells.Item(18, 3).Value = 1202.317307692308
$ws.Cells.Item(18, 4).Value = 1211.682692307692
$ws.Cells.Item(18, 5).Value = 1153.461538461539
$ws.Cells.Item(18, 6).Value = 1073.076923076923
$ws.Cells.Item(18, 7).Value = 1068
$ws.Cells.Item(19, 1).Value = 1164.894230769231
$ws.Cells.Item(19, 2).Value = 1168
$ws.Cells.Item(19, 3).Value = 1169.317307692308
$ws.Cells.Item(19, 4).Value = 1180.846153846154
$ws.Cells.Item(19, 5).Value = 1120.846153846154
$ws.Cells.Item(19, 6).Value = 1072.221153846154
$ws.Cells.Item(19, 7).Value = 1058.865384615385
$ws.Cells.Item(20, 1).Value = 1142.096153846154
$ws.Cells.Item(20, 2).Value = 1142.028301886792
$ws.Cells.Item(20, 3).Value = 1143.903846153846
$ws.Cells.Item(20, 4).Value = 1150.009615384615
$ws.Cells.Item(20, 5).Value = 1091.509615384615
$ws.Cells.Item(20, 6).Value = 1066.307692307692
$ws.Cells.Item(20, 7).Value = 1056.192307692308
$ws.Cells.Item(21, 1).Value = 1122.875
$ws.Cells.Item(21, 2).Value = 1117.245283018868
$ws.Cells.Item(21, 3).Value = 1124.163461538461
$ws.Cells.Item(21, 4).Value = 1126.826923076923
$ws.Cells.Item(21, 5).Value = 1077.307692307692
$ws.Cells.Item(21, 6).Value = 1059.269230769231
$ws.Cells.Item(21, 7).Value = 1054.855769230769
$ws.Cells.Item(22, 1).Value = 1103.461538461539
$ws.Cells.Item(22, 2).Value = 1102.283018867925
$ws.Cells.Item(22, 3).Value = 1105.759615384615
$ws.Cells.Item(22, 4).Value = 1110.105769230769
$ws.Cells.Item(22, 5).Value = 1069.019230769231
$ws.Cells.Item(22, 6).Value = 1056.471153846154
$ws.Cells.Item(22, 7).Value = 1049.865384615385
$ws.Cells.Item(23, 1).Value = 1085.134615384615
$ws.Cells.Item(23, 2).Value = 1083.88679245283
$ws.Cells.Item(23, 3).Value = 1082.230769230769
$ws.Cells.Item(23, 4).Value = 1090.442307692308
$ws.Cells.Item(23, 5).Value = 1059.211538461539
$ws.Cells.Item(23, 6).Value = 1043.134615384615
$ws.Cells.Item(23, 7).Value = 1047.730769230769
$ws.Cells.Item(24, 1).Value = 1068.644230769231
$ws.Cells.Item(24, 2).Value = 1066.735849056604
$ws.Cells.Item(24, 3).Value = 1067.961538461539
$ws.Cells.Item(24, 4).Value = 1075.923076923077
$ws.Cells.Item(24, 5).Value = 1056.221153846154
$ws.Cells.Item(24, 6).Value = 1052.653846153846
$ws.Cells.Item(24, 7).Value = 1047.259615384615
$ws.Cells.Item(25, 1).Value = 1059.673076923077
$ws.Cells.Item(25, 2).Value = 1057.820754716981
$ws.Cells.Item(25, 3).Value = 1058.644230769231
$ws.Cells.Item(25, 4).Value = 1063.173076923077
$ws.Cells.Item(25, 5).Value = 1052.009615384615
$ws.Cells.Item(25, 6).Value = 1045.807692307692
$ws.Cells.Item(25, 7).Value = 1042.807692307692
